$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price cells keep plain-number-looking values as TEXT (matching original inlineStr formatting)
$textCells = @("D5","D6","D8","D10","D12","D13","D14","D17","D19","D23","D24","D26","D30","D31","D33","D35","D36","D37","D38","D41","D43","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values row by row
$ws.Range("D2").Value = '51.795.10'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '2.787.17'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '352.27'
$ws.Range("E5").Value = '  -1.52%  '

$ws.Range("D6").Value = '109.07'
$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("E7").Value = '  -2.20%  '

$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +1.62%  '

$ws.Range("D10").Value = '39.91'
$ws.Range("E10").Value = '  -0.22%  '

$ws.Range("E11").Value = '  +2.65%  '

$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").Value = '20.17'
$ws.Range("E12").Value = '  +3.55%  '

$ws.Range("B13").Value = 'Dogecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D13").Value = '0.0838'
$ws.Range("E13").Value = '  -1.84%  '

$ws.Range("D14").Value = '7.68'
$ws.Range("E14").Value = '  +1.57%  '

$ws.Range("D15").Value = '3.224.02'
$ws.Range("E15").Value = '  -0.05%  '

$ws.Range("D16").Value = '2.795.26'
$ws.Range("E16").Value = '  +0.92%  '

$ws.Range("D17").Value = '0.927'
$ws.Range("E17").Value = '  -1.72%  '

$ws.Range("D18").Value = '51.759.04'
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("D19").Value = '7.77'
$ws.Range("E19").Value = '  +4.98%  '

$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("E21").Value = '  +1.39%  '

$ws.Range("D22").Value = '0.0₃0965'
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("D23").Value = '69.90'
$ws.Range("E23").Value = '  -0.51%  '

$ws.Range("D24").Value = '267.07'
$ws.Range("E24").Value = '  -2.63%  '

$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("D26").Value = '26.17'
$ws.Range("E26").Value = '  -1.91%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("E28").Value = '  +12.00%  '

$ws.Range("E29").Value = '  +0.47%  '

$ws.Range("D30").Value = '36.87'
$ws.Range("E30").Value = '  +7.34%  '

$ws.Range("D31").Value = '2.23'
$ws.Range("E31").Value = '  +0.76%  '

$ws.Range("E32").Value = '  +8.51%  '

$ws.Range("D33").Value = '51.84'
$ws.Range("E33").Value = '  +0.56%  '

$ws.Range("E34").Value = '  -3.11%  '

$ws.Range("D35").Value = '5.53'
$ws.Range("E35").Value = '  +5.21%  '

$ws.Range("D36").Value = '0.0831'
$ws.Range("E36").Value = '  -1.41%  '

$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("D38").Value = '18.52'
$ws.Range("E38").Value = '  +3.04%  '

$ws.Range("E39").Value = '  -2.54%  '

$ws.Range("E40").Value = '  -1.48%  '

$ws.Range("D41").Value = '2.55'
$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("E42").Value = '  -0.45%  '

$ws.Range("D43").Value = '120.72'
$ws.Range("E43").Value = '  -0.87%  '

$ws.Range("E44").Value = '  +0.18%  '

$ws.Range("E45").Value = '  -2.93%  '

$ws.Range("D46").Value = '2.123.45'
$ws.Range("E46").Value = '  +2.36%  '

$ws.Range("E47").Value = '  +2.01%  '

$ws.Range("E48").Value = '  +6.72%  '

$ws.Range("E49").Value = '  -2.19%  '

$ws.Range("D50").Value = '5.43'
$ws.Range("E50").Value = '  -5.03%  '

$ws.Range("E51").Value = '  +9.26%  '
